$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point precision on A3 (recalculated serial date value)
$ws.Range("A3").Value = 45875.0835671875

# Append new row 4 with the same date/time number format as A2/A3
$ws.Range("A4").Value = 45875.12516701442
$ws.Range("A4").NumberFormat = $ws.Range("A3").NumberFormat

$ws.Range("B4").Value = 2025
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = 13.72
$ws.Range("E4").Value = 92.56
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 4.47
$ws.Range("H4").Value = "NW"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "03:00:14"
